$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of cell address -> new value, derived from the source diff.
$updates = @{
    "D2" = "64.652.96"
    "E2" = "  -1.66%  "
    "D3" = "3.119.64"
    "E3" = "  -7.77%  "
    "E4" = "  -0.05%  "
    "D5" = "565.69"
    "E5" = "  -2.42%  "
    "D6" = "168.22"
    "E6" = "  -5.25%  "
    "D7" = "0.604"
    "E7" = "  -2.14%  "
    "D8" = "0.997"
    "E8" = "  -0.28%  "
    "D9" = "3.117.93"
    "E9" = "  -7.72%  "
    "E10" = "  -5.94%  "
    "E11" = "  -6.12%  "
    "E12" = "  -5.12%  "
    "D13" = "3.655.88"
    "E13" = "  -7.90%  "
    "D14" = "0.135"
    "E14" = "  +0.76%  "
    "D15" = "26.54"
    "E15" = "  -7.66%  "
    "D16" = "64.498.42"
    "E16" = "  -2.13%  "
    "E17" = "  -6.22%  "
    "D18" = "3.119.47"
    "E18" = "  -7.51%  "
    "D19" = "5.64"
    "E19" = "  -3.39%  "
    "D20" = "12.59"
    "E20" = "  -7.88%  "
    "D21" = "353.48"
    "E21" = "  -3.02%  "
    "D22" = "7.15"
    "E22" = "  -4.50%  "
    "E23" = "  +0.50%  "
    "D24" = "68.56"
    "E24" = "  -5.47%  "
    "D25" = "3.277.06"
    "E25" = "  -7.56%  "
    "D26" = "0.489"
    "E26" = "  -7.09%  "
    "D27" = "0.0000113"
    "E27" = "  -7.65%  "
    "D28" = "9.52"
    "E28" = "  -2.01%  "
    "D29" = "0.174"
    "E29" = "  -2.04%  "
    "E30" = "  +0.14%  "
    "D31" = "0.999"
    "E31" = "  -0.08%  "
    "E32" = "  -4.62%  "
    "D33" = "21.60"
    "E33" = "  -6.10%  "
    "D34" = "5.22"
    "E34" = "  -8.66%  "
    "B35" = "Aptos"
    "C35" = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
    "D35" = "6.50"
    "E35" = "  -6.45%  "
    "B36" = "Fetch.AI"
    "C36" = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
    "D36" = "1.18"
    "E36" = "  -5.33%  "
    "D37" = "157.86"
    "E37" = "  -2.39%  "
    "D38" = "1.41"
    "E38" = "  -7.16%  "
    "D39" = "0.825"
    "E39" = "  -3.58%  "
    "D40" = "25.79"
    "E40" = "  -5.09%  "
    "E41" = "  -1.28%  "
    "D42" = "2.634.79"
    "E42" = "  -1.38%  "
    "D43" = "2.39"
    "E43" = "  -7.66%  "
    "D44" = "6.04"
    "E44" = "  -2.26%  "
    "B45" = "OKB"
    "C45" = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
    "D45" = "39.21"
    "E45" = "  -0.65%  "
    "B46" = "Filecoin"
    "C46" = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
    "D46" = "4.11"
    "E46" = "  -5.08%  "
    "D47" = "0.0643"
    "E47" = "  -5.21%  "
    "D48" = "23.60"
    "E48" = "  -3.23%  "
    "D49" = "313.99"
    "E49" = "  -5.02%  "
    "E50" = "  -4.97%  "
    "E51" = "  -1.90%  "
}

foreach ($addr in $updates.Keys) {
    $cell = $ws.Range($addr)
    # Force text interpretation so numeric-looking strings (e.g. "64.652.96",
    # "0.997", "26.54") are not silently coerced into Double values, which
    # would corrupt formatting (trailing zeros, multi-dot "thousands" groups).
    $cell.NumberFormat = "@"
    $cell.Value = $updates[$addr]
    # Restore the default "Normal" style so the cell style stays identical
    # to before the edit (only the text content should change).
    $cell.Style = "Normal"
}
